# Logged Week 15 and simulated Week 16
# Updates cumulative season stats on the "Rushing" and "Receiving" sheets.

$wb = $excel.ActiveWorkbook

# --- Rushing sheet ---
$wsRushing = $wb.Worksheets.Item("Rushing")

# Row 2: B.Roethlisberger -> D2 (2DATT), E2 (3DATT), F2 (RZATT)
$wsRushing.Range("D2").Value = 1
$wsRushing.Range("E2").Value = 4
$wsRushing.Range("F2").Value = 3

# Row 4: N.Harris -> C4 (1DATT), D4 (2DATT), F4 (RZATT)
$wsRushing.Range("C4").Value = 143
$wsRushing.Range("D4").Value = 87
$wsRushing.Range("F4").Value = 26

# Row 8: D.Johnson -> C8 (1DATT)
$wsRushing.Range("C8").Value = 2

# Row 9: C.Claypool -> C9 (1DATT)
$wsRushing.Range("C9").Value = 6

# --- Receiving sheet ---
$wsReceiving = $wb.Worksheets.Item("Receiving")

# Row 2: N.Harris -> C2 (Short Target), D2 (Short Comp)
$wsReceiving.Range("C2").Value = 77
$wsReceiving.Range("D2").Value = 59

# Row 7: D.Johnson -> C7 (Short Target), D7 (Short Comp), E7 (Deep Target)
$wsReceiving.Range("C7").Value = 106
$wsReceiving.Range("D7").Value = 73
$wsReceiving.Range("E7").Value = 34

# Row 8: C.Claypool -> C8 (Short Target), D8 (Short Comp), E8 (Deep Target)
$wsReceiving.Range("C8").Value = 58
$wsReceiving.Range("D8").Value = 36
$wsReceiving.Range("E8").Value = 27

# Row 9: J.Washington -> C9 (Short Target), D9 (Short Comp), E9 (Deep Target), G9 (RZ Target), H9 (RZ Comp)
$wsReceiving.Range("C9").Value = 31
$wsReceiving.Range("D9").Value = 21
$wsReceiving.Range("E9").Value = 11
$wsReceiving.Range("G9").Value = 8
$wsReceiving.Range("H9").Value = 5

# Row 10: R.McCloud -> C10 (Short Target), E10 (Deep Target), G10 (RZ Target)
$wsReceiving.Range("C10").Value = 33
$wsReceiving.Range("E10").Value = 5
$wsReceiving.Range("G10").Value = 3

# Row 13: P.Freiermuth -> C13 (Short Target), D13 (Short Comp)
$wsReceiving.Range("C13").Value = 58
$wsReceiving.Range("D13").Value = 46

# Row 14: Z.Gentry -> C14 (Short Target), D14 (Short Comp)
$wsReceiving.Range("C14").Value = 11
$wsReceiving.Range("D14").Value = 9
